# Applies the "report eklendi" edit: expands the test-results sheet
# (OrnekSayfasi) from 14 data rows to 62, rewriting columns A (test
# case name) and B (result) for the full run, then restores the
# selection/view state left by the author (scrolled to row 10, with
# A18:A19 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, "New Account Creation", "Passed"),
    @(2, "New Account Creation", "Passed"),
    @(3, "Bill Pay", "Fail"),
    @(4, "Bill Pay", "Fail"),
    @(5, "Login with username and password", "Fail"),
    @(6, "Login with username and password", "Fail"),
    @(7, "Loan Application", "Fail"),
    @(8, "Register with username and password", "Passed"),
    @(9, "Login with username and password", "Passed"),
    @(10, "Register with username and password", "Passed"),
    @(11, "New Account Creation", "Passed"),
    @(12, "Bill Pay", "Passed"),
    @(13, "Update Contact Info", "Passed"),
    @(14, "Inter-Account Funds", "Passed"),
    @(15, "Loan Application", "Passed"),
    @(16, "Login with username and password", "Passed"),
    @(17, "New Account Creation", "Passed"),
    @(18, "Bill Pay", "Fail"),
    @(19, "Inter-Account Funds", "Fail"),
    @(20, "Loan Application", "Fail"),
    @(21, "Login with username and password", "Fail"),
    @(22, "Register with username and password", "Fail"),
    @(23, "New Account Creation", "Passed"),
    @(24, "Bill Pay", "Fail"),
    @(25, "Login with username and password", "Fail"),
    @(26, "Register with username and password", "Fail"),
    @(27, "New Account Creation", "Passed"),
    @(28, "Bill Pay", "Fail"),
    @(29, "Login with username and password", "Fail"),
    @(30, "Register with username and password", "Fail"),
    @(31, "New Account Creation", "Passed"),
    @(32, "Bill Pay", "Passed"),
    @(33, "Login with username and password", "Passed"),
    @(34, "Register with username and password", "Fail"),
    @(35, "New Account Creation", "Passed"),
    @(36, "Bill Pay", "Passed"),
    @(37, "Login with username and password", "Passed"),
    @(38, "Register with username and password", "Fail"),
    @(39, "New Account Creation", "Passed"),
    @(40, "Bill Pay", "Passed"),
    @(41, "Login with username and password", "Passed"),
    @(42, "Register with username and password", "Fail"),
    @(43, "New Account Creation", "Passed"),
    @(44, "Bill Pay", "Passed"),
    @(45, "Login with username and password", "Passed"),
    @(46, "Register with username and password", "Fail"),
    @(47, "New Account Creation", "Passed"),
    @(48, "Bill Pay", "Passed"),
    @(49, "Inter-Account Funds", "Fail"),
    @(50, "Loan Application", "Passed"),
    @(51, "Login with username and password", "Passed"),
    @(52, "Register with username and password", "Fail"),
    @(53, "New Account Creation", "Fail"),
    @(54, "Bill Pay", "Fail"),
    @(55, "Inter-Account Funds", "Fail"),
    @(56, "Loan Application", "Fail"),
    @(57, "Login with username and password", "Fail"),
    @(58, "Register with username and password", "Passed"),
    @(59, "New Account Creation", "Passed"),
    @(60, "Bill Pay", "Passed"),
    @(61, "Login with username and password", "Passed"),
    @(62, "Register with username and password", "Fail")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}

# Restore the view/selection state captured in the saved workbook:
# scrolled so row 10 is at the top, with A18:A19 selected (active cell A19).
[void]$ws.Range("A18:A19").Select()
try {
    $excel.ActiveWindow.ScrollRow = 10
    $excel.ActiveWindow.ScrollColumn = 1
} catch { }
